# This workbook tracks "sending windows" (time ranges). Two of the time
# ranges were edited in place:
#   B9:  "18:55 - 18:59" -> "19:40 - 19:44"
#   B10: "19:00 - 19:04" -> "19:45 - 19:49"
# (Excel will naturally rebuild/compact the shared-strings table as a
# side effect of this edit, which accounts for the other apparent index
# churn in the diff even though the displayed text of those cells is
# unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "19:40 - 19:44"
$ws.Range("B10").Value = "19:45 - 19:49"

# Reflect the author's final view/selection state (scrolled so row 7 is at
# the top, with B13 as the active/selected cell).
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B13").Select()

$wb.Save()
